# Update "想去人数" (F column) values per latest scraped output
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 372
$ws.Range("F6").Value = 552
$ws.Range("F9").Value = 11812
$ws.Range("F13").Value = 2124
$ws.Range("F20").Value = 255
$ws.Range("F24").Value = 2406
$ws.Range("F26").Value = 3692
$ws.Range("F27").Value = 3692
$ws.Range("F28").Value = 1078
$ws.Range("F33").Value = 1000
$ws.Range("F40").Value = 3692
$ws.Range("F41").Value = 4457
$ws.Range("F42").Value = 5494
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 4160
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F8").Value = 372
$ws.Range("F9").Value = 552
$ws.Range("F11").Value = 11812
$ws.Range("F15").Value = 2124
$ws.Range("F20").Value = 255
$ws.Range("F21").Value = 4160
$ws.Range("F24").Value = 3692
$ws.Range("F25").Value = 1078
$ws.Range("F30").Value = 1000
$ws.Range("F36").Value = 4457
